$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.246.45"
$ws.Range("E2").Value = "  +0.19%  "

$ws.Range("D3").Value = "3.099.52"
$ws.Range("E3").Value = "  -0.98%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'580.89"
$ws.Range("E5").Value = "  +0.26%  "

$ws.Range("D6").Value = "'170.84"
$ws.Range("E6").Value = "  -1.80%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "3.097.63"
$ws.Range("E8").Value = "  -0.96%  "

$ws.Range("D9").Value = "'0.517"
$ws.Range("E9").Value = "  -1.29%  "

$ws.Range("D10").Value = "'6.44"
$ws.Range("E10").Value = "  -0.53%  "

$ws.Range("E11").Value = "  -2.33%  "

$ws.Range("D12").Value = "'0.475"
$ws.Range("E12").Value = "  -1.12%  "

$ws.Range("D13").Value = "'0.0000244"
$ws.Range("E13").Value = "  -2.17%  "

$ws.Range("D14").Value = "'36.51"
$ws.Range("E14").Value = "  -2.08%  "

$ws.Range("D15").Value = "'0.121"
$ws.Range("E15").Value = "  -1.91%  "

$ws.Range("D16").Value = "3.615.86"
$ws.Range("E16").Value = "  -0.90%  "

$ws.Range("D17").Value = "67.202.71"
$ws.Range("E17").Value = "  +0.17%  "

$ws.Range("D18").Value = "'7.07"
$ws.Range("E18").Value = "  -1.47%  "

$ws.Range("B19").Value = "WrappedEther"
$ws.Range("C19").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D19").Value = "3.100.43"
$ws.Range("E19").Value = "  -1.06%  "

$ws.Range("B20").Value = "Chainlink"
$ws.Range("C20").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D20").Value = "'16.56"
$ws.Range("E20").Value = "  +2.58%  "

$ws.Range("D21").Value = "'486.83"
$ws.Range("E21").Value = "  +0.35%  "

$ws.Range("D22").Value = "'7.78"
$ws.Range("E22").Value = "  +1.30%  "

$ws.Range("D23").Value = "'0.694"
$ws.Range("E23").Value = "  -3.00%  "

$ws.Range("D24").Value = "'83.52"
$ws.Range("E24").Value = "  -0.73%  "

$ws.Range("D25").Value = "'12.99"
$ws.Range("E25").Value = "  -2.65%  "

$ws.Range("D26").Value = "'2.27"
$ws.Range("E26").Value = "  -2.40%  "

$ws.Range("D27").Value = "'10.46"
$ws.Range("E27").Value = "  +4.20%  "

$ws.Range("E28").Value = "  -0.01%  "

$ws.Range("D29").Value = "'7.74"
$ws.Range("E29").Value = "  -3.14%  "

$ws.Range("E30").Value = "  -3.35%  "

$ws.Range("E31").Value = "  -1.16%  "

$ws.Range("D32").Value = "'28.17"
$ws.Range("E32").Value = "  -2.30%  "

$ws.Range("E33").Value = "  -1.63%  "

$ws.Range("D34").Value = "0.0₃0933"
$ws.Range("E34").Value = "  -6.05%  "

$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("D36").Value = "'5.76"
$ws.Range("E36").Value = "  -2.39%  "

$ws.Range("D37").Value = "'0.963"
$ws.Range("E37").Value = "  -2.04%  "

$ws.Range("D38").Value = "'46.33"
$ws.Range("E38").Value = "  -3.01%  "

$ws.Range("D39").Value = "'2.01"
$ws.Range("E39").Value = "  -4.44%  "

$ws.Range("D40").Value = "'0.123"
$ws.Range("E40").Value = "  +0.81%  "

$ws.Range("D41").Value = "'0.304"
$ws.Range("E41").Value = "  -2.27%  "

$ws.Range("D42").Value = "'8.39"
$ws.Range("E42").Value = "  -3.09%  "

$ws.Range("D43").Value = "2.790.65"
$ws.Range("E43").Value = "  -1.99%  "

$ws.Range("D44").Value = "'379.98"
$ws.Range("E44").Value = "  -0.65%  "

$ws.Range("E45").Value = "  -6.68%  "

$ws.Range("D46").Value = "'0.0348"
$ws.Range("E46").Value = "  -2.93%  "

$ws.Range("D47").Value = "'135.20"
$ws.Range("E47").Value = "  -0.20%  "

$ws.Range("E48").Value = "  -0.01%  "

$ws.Range("D49").Value = "'24.71"
$ws.Range("E49").Value = "  -0.65%  "

$ws.Range("D50").Value = "'2.17"
$ws.Range("E50").Value = "  -2.14%  "

$ws.Range("E51").Value = "  -1.97%  "
